$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.210.80'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.91%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''3.498.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +0.76%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.09%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''602.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.48%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''143.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -1.94%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''3.498.66'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.73%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.15%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.474'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -0.68%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''8.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +7.54%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.136'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -3.90%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.412'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -2.15%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''4.057.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +0.04%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''0.0000203'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -3.93%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''30.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -3.28%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''3.481.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +0.26%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = '''WrappedBTC'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = '''66.146.60'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -1.03%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = '''TRON'
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = '''https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = '''0.116'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.49%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''10.52'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +5.54%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''6.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -3.58%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''14.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -2.95%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''421.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -2.53%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''0.586'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -2.98%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''77.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -1.54%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.04%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''0.0000117'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -2.41%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''9.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -3.23%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''8.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -3.96%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''2.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.75%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.31%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.163'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -2.17%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''1.48'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -5.79%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''25.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.47%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''3.476.32'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +0.47%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -0.07%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''1.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -3.79%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''5.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -5.38%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''7.66'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -2.33%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  -0.03%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''169.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -2.80%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.0871'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.70%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.894'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.32%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''5.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -4.60%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = '''  -8.69%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''45.72'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -1.15%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''26.45'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -7.35%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''1.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -1.01%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''dogwifhat'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''2.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -2.36%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = '''Cosmos'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''7.13'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -3.74%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.937'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -3.52%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.236'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -3.16%  '
$ws.Range('E51').Style = 'Normal'
